$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

# Row 2 - Bitcoin
Set-TextValue "D2" "58.666.31"
$ws.Range("E2").Value = "  -1.32%  "

# Row 3 - Ethereum
Set-TextValue "D3" "2.631.78"
$ws.Range("E3").Value = "  +0.97%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.01%  "

# Row 5 - BNB
Set-TextValue "D5" "536.56"
$ws.Range("E5").Value = "  -0.10%  "

# Row 6 - Solana
Set-TextValue "D6" "143.34"
$ws.Range("E6").Value = "  +1.34%  "

# Row 7 - USDC
Set-TextValue "D7" "0.999"
$ws.Range("E7").Value = "  -0.09%  "

# Row 8 - XRP
$ws.Range("E8").Value = "  -0.02%  "

# Row 9 - LidoStakedEther
Set-TextValue "D9" "2.641.36"
$ws.Range("E9").Value = "  +1.10%  "

# Row 10 - Toncoin
Set-TextValue "D10" "7.03"
$ws.Range("E10").Value = "  +7.97%  "

# Row 11 - Dogecoin
$ws.Range("E11").Value = "  -1.62%  "

# Row 12 - Cardano
$ws.Range("E12").Value = "  -0.03%  "

# Row 13 - TRON
$ws.Range("E13").Value = "  +0.96%  "

# Row 14 - WrappedliquidstakedEther2.0
Set-TextValue "D14" "3.095.92"

# Row 15 - WrappedBTC
Set-TextValue "D15" "58.611.27"
$ws.Range("E15").Value = "  -1.28%  "

# Row 16 - Avalanche
$ws.Range("E16").Value = "  +1.12%  "

# Row 17 - WrappedEther
Set-TextValue "D17" "2.645.10"
$ws.Range("E17").Value = "  +0.22%  "

# Row 18 - ShibaInu
$ws.Range("E18").Value = "  -0.78%  "

# Row 19 - Polkadot
$ws.Range("E19").Value = "  +1.17%  "

# Row 20 - BitcoinCash
Set-TextValue "D20" "334.34"
$ws.Range("E20").Value = "  -2.09%  "

# Row 21 - Chainlink
$ws.Range("E21").Value = "  +0.65%  "

# Row 22 - Uniswap
$ws.Range("E22").Value = "  -1.84%  "

# Row 23 - Dai
$ws.Range("E23").Value = "  -0.02%  "

# Row 24 - Litecoin
Set-TextValue "D24" "66.26"
$ws.Range("E24").Value = "  -1.82%  "

# Row 25 - Polygon
Set-TextValue "D25" "0.416"
$ws.Range("E25").Value = "  +1.66%  "

# Row 26 - Kaspa
$ws.Range("E26").Value = "  -0.69%  "

# Row 27 - Binance-PegBSC-USD
$ws.Range("E27").Value = "  +0.41%  "

# Row 28 - InternetComputer(DFINITY)
Set-TextValue "D28" "7.17"
$ws.Range("E28").Value = "  -0.85%  "

# Row 29 - PEPE
$ws.Range("E29").Value = "  -0.73%  "

# Row 31 - PancakeSwap
$ws.Range("E31").Value = "  -1.34%  "

# Row 32 - Aptos
Set-TextValue "D32" "5.86"
$ws.Range("E32").Value = "  +0.79%  "

# Row 33 - EthereumClassic
$ws.Range("E33").Value = "  -0.52%  "

# Row 34 - Monero
Set-TextValue "D34" "150.81"
$ws.Range("E34").Value = "  +0.65%  "

# Row 35 - NEARProtocol
$ws.Range("E35").Value = "  -1.49%  "

# Row 36 - OKB
Set-TextValue "D36" "37.17"
$ws.Range("E36").Value = "  +0.03%  "

# Row 37 - ImmutableX
$ws.Range("E37").Value = "  -0.63%  "

# Row 38 - SuiNetwork
Set-TextValue "D38" "0.850"
$ws.Range("E38").Value = "  +2.03%  "

# Row 39 - Stacks
$ws.Range("E39").Value = "  -2.85%  "

# Row 40 - Fetch.AI
$ws.Range("E40").Value = "  -1.29%  "

# Row 41 - Filecoin
$ws.Range("E41").Value = "  +1.37%  "

# Row 42 - Bittensor
Set-TextValue "D42" "282.74"
$ws.Range("E42").Value = "  +3.36%  "

# Row 43 - FirstDigitalUSD
Set-TextValue "D43" "0.999"
$ws.Range("E43").Value = "  -0.13%  "

# Row 44 - Mantle
$ws.Range("E44").Value = "  +0.33%  "

# Row 45 - WhiteBITCoin
Set-TextValue "D45" "10.71"
$ws.Range("E45").Value = "  -0.23%  "

# Row 46 - was EnergySwap, now Hedera
$ws.Range("B46").Value = "Hedera"
$ws.Range("C46").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue "D46" "0.0533"
$ws.Range("E46").Value = "  +1.81%  "

# Row 47 - was Hedera, now EnergySwap
$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue "D47" "19.10"
$ws.Range("E47").Value = "  +3.19%  "

# Row 48 - Stellar
$ws.Range("E48").Value = "  -2.03%  "

# Row 49 - VeChain
$ws.Range("E49").Value = "  +0.98%  "

# Row 50 - Maker
Set-TextValue "D50" "1.949.87"
$ws.Range("E50").Value = "  +0.03%  "

# Row 51 - RenderToken
Set-TextValue "D51" "4.47"
$ws.Range("E51").Value = "  -0.95%  "
